# Update countries & provincias Spain
# Applies the 20-May-2020 19:35 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp
#  - re-ranks a few countries whose case counts caused them to swap/move
#    position in the (descending, by total cases) table
#  - refreshes the numeric columns (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#    for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: "Datos actualizados ..." timestamp ----
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 19:35"

# ---- Row 9: Italia (unchanged name, refreshed numbers) ----
$ws.Cells.Item(9, 2).Value = 227364
$ws.Cells.Item(9, 3).Value = 665
$ws.Cells.Item(9, 4).Value = 132282
$ws.Cells.Item(9, 5).Value = 62752
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 161
$ws.Cells.Item(9, 8).Value = 32330

# ---- Row 14: India (unchanged name, refreshed numbers) ----
$ws.Cells.Item(14, 2).Value = 111750
$ws.Cells.Item(14, 3).Value = 5275
$ws.Cells.Item(14, 4).Value = 45422
$ws.Cells.Item(14, 5).Value = 62894
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 132
$ws.Cells.Item(14, 8).Value = 3434

# ---- Rows 57-58: Australia / Marruecos swap ranking ----
$ws.Cells.Item(57, 1).Value = "Marruecos"
$ws.Cells.Item(57, 2).Value = 7133
$ws.Cells.Item(57, 3).Value = 110
$ws.Cells.Item(57, 4).Value = 4098
$ws.Cells.Item(57, 5).Value = 2841
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 194

$ws.Cells.Item(58, 1).Value = "Australia"
$ws.Cells.Item(58, 2).Value = 7079
$ws.Cells.Item(58, 3).Value = 11
$ws.Cells.Item(58, 4).Value = 6444
$ws.Cells.Item(58, 5).Value = 535
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 100

# ---- Rows 114-115: Zambia / Paraguay swap ranking ----
$ws.Cells.Item(114, 1).Value = "Zambia"
$ws.Cells.Item(114, 2).Value = 833
$ws.Cells.Item(114, 3).Value = 4
$ws.Cells.Item(114, 4).Value = 242
$ws.Cells.Item(114, 5).Value = 580
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 11

$ws.Cells.Item(115, 1).Value = "Paraguay"
$ws.Cells.Item(115, 2).Value = 832
$ws.Cells.Item(115, 3).Value = 60
$ws.Cells.Item(115, 4).Value = 197
$ws.Cells.Item(115, 5).Value = 628
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 7

# ---- Rows 127-135: "Estado de Palestina" moves up in ranking, ----
# ---- shifting Sierra Leona, Republica del Chad, Jamaica, Tanzania, ----
# ---- Reunion, Taiwan, Nepal and Congo down by one position, each ----
# ---- with refreshed numbers. ----
$ws.Cells.Item(127, 1).Value = "Estado de Palestina"
$ws.Cells.Item(127, 2).Value = 577
$ws.Cells.Item(127, 3).Value = 186
$ws.Cells.Item(127, 4).Value = 346
$ws.Cells.Item(127, 5).Value = 229
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 2

$ws.Cells.Item(128, 1).Value = "Sierra Leona"
$ws.Cells.Item(128, 2).Value = 570
$ws.Cells.Item(128, 3).Value = 36
$ws.Cells.Item(128, 4).Value = 205
$ws.Cells.Item(128, 5).Value = 331
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 34

$ws.Cells.Item(129, 1).Value = "Republica del Chad"
$ws.Cells.Item(129, 2).Value = 565
$ws.Cells.Item(129, 3).Value = 20
$ws.Cells.Item(129, 4).Value = 177
$ws.Cells.Item(129, 5).Value = 331
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 57

$ws.Cells.Item(130, 1).Value = "Jamaica"
$ws.Cells.Item(130, 2).Value = 520
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 145
$ws.Cells.Item(130, 5).Value = 366
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 9

$ws.Cells.Item(131, 1).Value = "Tanzania"
$ws.Cells.Item(131, 2).Value = 509
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 183
$ws.Cells.Item(131, 5).Value = 305
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 21

$ws.Cells.Item(132, 1).Value = "Reunion"
$ws.Cells.Item(132, 2).Value = 447
$ws.Cells.Item(132, 3).Value = 1
$ws.Cells.Item(132, 4).Value = 411
$ws.Cells.Item(132, 5).Value = 35
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 1
$ws.Cells.Item(132, 8).Value = 1

$ws.Cells.Item(133, 1).Value = "Taiwan"
$ws.Cells.Item(133, 2).Value = 440
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 402
$ws.Cells.Item(133, 5).Value = 31
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 7

$ws.Cells.Item(134, 1).Value = "Nepal"
$ws.Cells.Item(134, 2).Value = 427
$ws.Cells.Item(134, 3).Value = 25
$ws.Cells.Item(134, 4).Value = 45
$ws.Cells.Item(134, 5).Value = 380
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 2

$ws.Cells.Item(135, 1).Value = "Congo"
$ws.Cells.Item(135, 2).Value = 420
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 132
$ws.Cells.Item(135, 5).Value = 273
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 15

# ---- Row 145: Ruanda (unchanged name, refreshed numbers) ----
$ws.Cells.Item(145, 2).Value = 314
$ws.Cells.Item(145, 3).Value = 6
$ws.Cells.Item(145, 4).Value = 216
$ws.Cells.Item(145, 5).Value = 98

# ---- Row 169: Aruba (unchanged name, refreshed numbers) ----
$ws.Cells.Item(169, 4).Value = 95
$ws.Cells.Item(169, 5).Value = 3
